# "updated Medical Decision problem"
# Update parameter values on Sheet1/Sheet2/Sheet3 and move the active
# selection/tab from Sheet1 to Sheet3.

$wb = $excel.ActiveWorkbook

# --- Sheet1: high_low policy params ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("B3").Value = 6      # high_low param1: 5 -> 6
$ws1.Range("C3").Value = 13     # high_low param2: 15 -> 13
$ws1.Range("B4").Value = 1.9    # track param1: 1 -> 1.9

# --- Sheet2: grid-search bounds ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A2").Value = 0      # low_min: 4 -> 0
$ws2.Range("B2").Value = 0      # low_max: 7 -> 0
$ws2.Range("C2").Value = 2      # high_min: 15 -> 2
$ws2.Range("F2").Value = 2      # track_max: 4 -> 2
$ws2.Range("G2").Value = 0.1    # increment_size: 1 -> 0.1

# --- Sheet3: evaluation parameters (Medical Decision problem) ---
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("C2").Value = 20     # TimeHorizon: 100 -> 20
$ws3.Range("I2").Value = 200    # Iterations: 10 -> 200

# --- Restore each sheet's last selection ---
$ws1.Activate()
$ws1.Range("D14").Select()

$ws2.Activate()
$ws2.Range("G3").Select()

# Sheet3 ends up the active/selected tab (activeTab 1 -> 3)
$ws3.Activate()
$ws3.Range("B8").Select()
